$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 124, pushing the existing rows 124-137 down to 125-138
$ws.Rows.Item(124).Insert()

# Populate the newly inserted row 124 with the new weekly record
$ws.Cells.Item(124, 1).Value = 3
$ws.Cells.Item(124, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(124, 3).Value = "Coquimbo"
$ws.Cells.Item(124, 4).Value = 44449
$ws.Cells.Item(124, 5).Value = 5
$ws.Cells.Item(124, 6).Value = 100112001
$ws.Cells.Item(124, 7).Value = "Berenjena"
$ws.Cells.Item(124, 8).Value = "Sin especificar"
$ws.Cells.Item(124, 9).Value = "Primera"
$ws.Cells.Item(124, 10).Value = 85
$ws.Cells.Item(124, 11).Value = 10000
$ws.Cells.Item(124, 12).Value = 10500
$ws.Cells.Item(124, 13).Value = 10235
$ws.Cells.Item(124, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(124, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(124, 16).Value = 171
$ws.Cells.Item(124, 17).Value = 60
$ws.Cells.Item(124, 18).Value = "Hortaliza"
